# Determine whether expression between genes within clusters is more similar
# than between all pairwise comparisons.
#
# The "Clustering" column (K) on Sheet1 previously only labelled rows as
# "single", "cluster" or "double". Every distinct named gene cluster that
# was tagged "cluster"/"double" now gets its own unique label ("cluster 1",
# "cluster 2", ... "cluster 33") so that downstream analysis can group genes
# by their specific cluster rather than lumping every clustered gene
# together.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry is the K-column range belonging to one distinct gene cluster,
# paired with the new unique cluster label that replaces the generic
# "cluster"/"double" text.
$clusterRanges = @(
    @("K3:K6",    "cluster 1"),
    @("K8:K10",   "cluster 2"),
    @("K11:K12",  "cluster 3"),
    @("K19:K21",  "cluster 4"),
    @("K23:K30",  "cluster 5"),
    @("K33:K36",  "cluster 6"),
    @("K38:K39",  "cluster 7"),
    @("K40:K46",  "cluster 8"),
    @("K50:K61",  "cluster 9"),
    @("K64:K65",  "cluster 10"),
    @("K66:K67",  "cluster 11"),
    @("K69:K70",  "cluster 12"),
    @("K73:K76",  "cluster 13"),
    @("K77:K79",  "cluster 14"),
    @("K81:K82",  "cluster 15"),
    @("K84:K86",  "cluster 16"),
    @("K89:K90",  "cluster 17"),
    @("K92:K100", "cluster 18"),
    @("K101:K104","cluster 19"),
    @("K108:K109","cluster 20"),
    @("K111:K112","cluster 21"),
    @("K117:K119","cluster 22"),
    @("K120:K122","cluster 23"),
    @("K131:K132","cluster 24"),
    @("K133:K137","cluster 25"),
    @("K139:K141","cluster 26"),
    @("K143:K146","cluster 27"),
    @("K147:K153","cluster 28"),
    @("K154:K155","cluster 29"),
    @("K156:K162","cluster 30"),
    @("K163:K165","cluster 31"),
    @("K167:K169","cluster 32"),
    @("K174:K175","cluster 33")
)

foreach ($entry in $clusterRanges) {
    $addr = $entry[0]
    $label = $entry[1]
    $ws.Range($addr).Value = $label
}

# Restore the sheet view scroll position / active cell saved with the
# workbook (topLeftCell moved from A135 to A136, selection moved from A174
# to K175).
$ws.Activate()
$ws.Range("K175").Select()
$excel.ActiveWindow.ScrollRow = 136
$excel.ActiveWindow.ScrollColumn = 1
